$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti)
$newRows = @(
    @(44326, 1, 10, 64.58696635019054),
    @(44327, 1, 11, 71.04566298520959),
    @(44328, 0, 9, 58.12826971517148),
    @(44329, 1, 8, 51.66957308015243)
)

$lastRow = 251

foreach ($rowData in $newRows) {
    $newRowIndex = $lastRow + 1

    # Copy formatting from the cell above (column A carries the date style)
    $ws.Cells.Item($lastRow, 1).Copy()
    $ws.Cells.Item($newRowIndex, 1).PasteSpecial(-4122)

    $ws.Cells.Item($newRowIndex, 1).Value = $rowData[0]
    $ws.Cells.Item($newRowIndex, 2).Value = $rowData[1]
    $ws.Cells.Item($newRowIndex, 3).Value = $rowData[2]
    $ws.Cells.Item($newRowIndex, 4).Value = $rowData[3]

    $lastRow = $newRowIndex
}

Write-Output "Appended $($newRows.Count) rows, last row is now $lastRow"
